# HOTFIX: Use BackgroundTasks for Workday Email to prevent 502 Timeouts
#
# Fix the typo "Juan Parez" -> "Juan Perez" in I2, replace K2's value
# ("CIAL_ALIMENTOS") with a test value ("HOLA prueba"), and remove the
# extra sample rows (3, 4, 5) that were left over from testing, keeping
# only the header row and a single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unneeded sample rows 3-5 (keep header row 1 and data row 2)
$ws.Range("A3:K5").EntireRow.Delete()

# Fix values in the remaining data row
$ws.Range("I2").Value = "Juan Perez"
$ws.Range("K2").Value = "HOLA prueba"

# Update the active selection to match the saved workbook state
$ws.Range("F12").Select()
